$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "218.13"). Force them to remain plain text, matching the source data,
# by briefly marking the cell as Text before assignment, then clearing the
# number-format override again so no stray style is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "26.107.19"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.651.07"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "218.13"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "0.06280"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "20.45"
$ws.Range("D11").Value = "0.07797"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "4.476"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("D13").Value = "1.657.98"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "1.878.58"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "0.5519"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "0.0₅7995"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("D17").Value = "64.70"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "26.094.72"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "4.617"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "194.08"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").Value = "5.937"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").Value = "146.60"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "0.1201"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("D27").Value = "7.149"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "15.87"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").Value = "1.481"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "0.05688"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").Value = "1.267"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "3.471"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("D33").Value = "3.331"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("D35").Value = "2.797"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "0.9466"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.414"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "0.5652"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").Value = "0.01588"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("D40").Value = "5.914"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "1.058.86"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "0.8411"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("D44").Value = "103.24"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "1.789.38"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").Value = "0.05416"
$ws.Range("E48").Value = "  +4.95%  "
$ws.Range("D49").Value = "1.008"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "7.942"
$ws.Range("E51").Value = "  -1.38%  "

# Remove the temporary text formatting now that the literal values are stored,
# so cell styling matches the original (unstyled) cells again.
$ws.Range("D5").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D51").ClearFormats()
